$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# Remove the old row 3 (ID 133 / Accrual on 41671) - the remaining rows
# shift up, bringing the old rows 4-7 into rows 3-6 unchanged.
$ws.Range("A3").EntireRow.Delete()

# Renumber the ID column (A) for the resulting 5 data rows.
$ws.Range("A2").Value = 3502
$ws.Range("A3").Value = 3501
$ws.Range("A4").Value = 3500
$ws.Range("A5").Value = 3497
$ws.Range("A6").Value = 3496

# Update the active selection to match the saved view state.
$ws.Range("D4").Select()
